$wb = $excel.ActiveWorkbook

# Overview sheet: update "Latest HO Xliff Generate Date" (G2)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-23 13:07:39"

# zh-cn sheet: update "Correspond Handoff Datetime" (H2) and
# "Correspond Handback DateTime" (K2)
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-23 13:07:34"
$wsZhCn.Range("K2").Value = "2016-08-23 13:07:52"

# de-de sheet: "Correspond Handoff Datetime" (H2) shares the same
# underlying shared-string slot as Overview!G2 ("Latest HO Xliff
# Generate Date"), so it picks up the same new value; also update
# "Correspond Handback DateTime" (K2)
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-23 13:07:39"
$wsDeDe.Range("K2").Value = "2016-08-23 13:07:59"
